# Weekly update: insert a new price-report row for "Berenjena" (Terminal La
# Palmera de La Serena) at row 188, pushing all existing data rows (188-251)
# down by one (to 189-252), and fill in the new row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 188; this shifts rows
# 188..251 down to 189..252 and extends the used range accordingly.
$ws.Rows.Item(188).Insert()

# Populate the newly inserted row 188 with this week's data.
$ws.Cells.Item(188, 1).Value  = 8
$ws.Cells.Item(188, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(188, 3).Value  = "Coquimbo"
$ws.Cells.Item(188, 4).Value  = 45093
$ws.Cells.Item(188, 5).Value  = 4
$ws.Cells.Item(188, 6).Value  = 100112001
$ws.Cells.Item(188, 7).Value  = "Berenjena"
$ws.Cells.Item(188, 8).Value  = "Sin especificar"
$ws.Cells.Item(188, 9).Value  = "Primera"
$ws.Cells.Item(188, 10).Value = 300
$ws.Cells.Item(188, 11).Value = 9000
$ws.Cells.Item(188, 12).Value = 10000
$ws.Cells.Item(188, 13).Value = 9500
$ws.Cells.Item(188, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(188, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(188, 16).Value = 190
$ws.Cells.Item(188, 17).Value = 50
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Apply the same date number format used by the other rows in column D.
$ws.Cells.Item(188, 4).NumberFormat = $ws.Cells.Item(189, 4).NumberFormat
